$d = $word.ActiveDocument

# --- Remove the last three rows (Cyber/6.4, Nozomi/5.2, Inizio/6.2) ---
$t = $d.Tables.Item(1)
$t.Rows.Item(6).Delete()
$t.Rows.Item(5).Delete()
$t.Rows.Item(4).Delete()

# --- Update remaining participant names / scores ---
$d.Content.Find.Execute("Altan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Yamai", 2)
$d.Content.Find.Execute("4.8", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4.7", 2)
$d.Content.Find.Execute("Goromi", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Kalani", 2)
$d.Content.Find.Execute("3.2", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5.7", 2)
